$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23/24: fill in the new "name" / "RoleName" column-B header cells ---
$ws.Range("B23").Value = "name"
$ws.Range("B24").Value = "RoleName"

# --- Rows 25-28: mirror column A into the (until now empty) column B ---
$ws.Range("B25").Value = "Tutor"
$ws.Range("B26").Value = "Student"
$ws.Range("B27").Value = "GradStudent"
$ws.Range("B28").Value = "AccountManager"

# --- Rows 29-34: mirror column A into column B, but these rows use the
#     "Menlo" cell style (s=12) instead of column B's default style, so copy
#     the formatting from column A before writing the values. ---
$ws.Range("A29").Copy()
$ws.Range("B29:B34").PasteSpecial(-4122)

$ws.Range("B29").Value = "SystemAdmin"
$ws.Range("B30").Value = "Administrator"
$ws.Range("B31").Value = "Anonymous"
$ws.Range("B32").Value = "ExecEngine"
$ws.Range("B33").Value = "Janitor"
$ws.Range("B34").Value = "SYSTEM"

# --- Insert three new rows (35-37) for User / FormalAmpersand / Type_32_Checker,
#     pushing the old rows 37-40 ([PF_NavMenuItem]..Logout) down to 40-43. ---
$ws.Rows("35:37").Insert()

# New rows inherit row 34's formatting on A/C already; copy it onto B too.
$ws.Range("A35").Copy()
$ws.Range("B35:B37").PasteSpecial(-4122)

$ws.Range("A35").Value = "User"
$ws.Range("B35").Value = "User"

$ws.Range("A36").Value = "FormalAmpersand"
$ws.Range("B36").Value = "FormalAmpersand"

$ws.Range("A37").Value = "Type_32_Checker"
$ws.Range("B37").Value = "Type_32_Checker"

# --- Update the sheet selection to match the authored state ---
$ws.Range("A38").Select()
